# Auto-generated edit script: updates cached market-price/profit
# figures on the Leve profit sheets (scheduled-runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1275.9651
$ws.Range("I15").Value = 1275.9651
$ws.Range("K15").Value = 3827.8953
$ws.Range("M15").Value = -3658.8953

$ws.Range("H76").Value = 45836332
$ws.Range("I76").Value = 45836332
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 45836332
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -45836017
$ws.Range("N76").Value = $null

$ws.Range("H79").Value = 45836332
$ws.Range("I79").Value = 45836332
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 45836332
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -45835240
$ws.Range("N79").Value = $null

$ws.Range("H115").Value = 601
$ws.Range("I115").Value = 601
$ws.Range("K115").Value = 1803
$ws.Range("M115").Value = -236

$ws.Range("H132").Value = 2407.4167
$ws.Range("I132").Value = 2355.5652
$ws.Range("K132").Value = 7066.6956
$ws.Range("M132").Value = -4536.6956

$ws.Range("H137").Value = 8198167.5
$ws.Range("I137").Value = 1474.1177
$ws.Range("J137").Value = 18519930
$ws.Range("K137").Value = 4422.3531
$ws.Range("L137").Value = 55559790
$ws.Range("M137").Value = -1872.3531
$ws.Range("N137").Value = -55564890

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34421.062
$ws.Range("I32").Value = 30609.318
$ws.Range("J32").Value = 90326.664
$ws.Range("K32").Value = 30609.318
$ws.Range("L32").Value = 90326.664
$ws.Range("M32").Value = -30322.318
$ws.Range("N32").Value = -90900.664

$ws.Range("H110").Value = 491.3143
$ws.Range("I110").Value = 470.51724
$ws.Range("J110").Value = 591.8333
$ws.Range("K110").Value = 470.51724
$ws.Range("L110").Value = 591.8333
$ws.Range("M110").Value = 1574.48276
$ws.Range("N110").Value = -4681.8333

$ws.Range("H132").Value = 1708.2291
$ws.Range("I132").Value = 1237.3823
$ws.Range("J132").Value = 2851.7144
$ws.Range("K132").Value = 3712.1469
$ws.Range("L132").Value = 8555.143199999999
$ws.Range("M132").Value = -1182.1469
$ws.Range("N132").Value = -13615.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 47651104
$ws.Range("I20").Value = 130557
$ws.Range("J20").Value = 62501276
$ws.Range("K20").Value = 130557
$ws.Range("L20").Value = 62501276
$ws.Range("M20").Value = -130310
$ws.Range("N20").Value = -62501770

$ws.Range("H94").Value = 1091.6364
$ws.Range("I94").Value = 972.44446
$ws.Range("K94").Value = 972.44446
$ws.Range("M94").Value = -521.44446

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1356.8715
$ws.Range("I68").Value = 731.1515000000001
$ws.Range("J68").Value = 1914.9459
$ws.Range("K68").Value = 2193.4545
$ws.Range("L68").Value = 5744.8377
$ws.Range("M68").Value = -1382.4545
$ws.Range("N68").Value = -7366.8377

$ws.Range("H71").Value = 1356.8715
$ws.Range("I71").Value = 731.1515000000001
$ws.Range("J71").Value = 1914.9459
$ws.Range("K71").Value = 6580.3635
$ws.Range("L71").Value = 17234.5131
$ws.Range("M71").Value = -2524.3635
$ws.Range("N71").Value = -25346.5131

$ws.Range("H113").Value = 433.02
$ws.Range("I113").Value = 386.57895
$ws.Range("J113").Value = 443.91357
$ws.Range("K113").Value = 1159.73685
$ws.Range("L113").Value = 1331.74071
$ws.Range("M113").Value = 1010.26315
$ws.Range("N113").Value = -5671.74071

$ws.Range("H122").Value = 499.21054
$ws.Range("I122").Value = 417.875
$ws.Range("K122").Value = 3760.875
$ws.Range("M122").Value = -1310.875

$ws.Range("H131").Value = 19585.018
$ws.Range("I131").Value = 83912.5
$ws.Range("J131").Value = 2041.159
$ws.Range("K131").Value = 251737.5
$ws.Range("L131").Value = 6123.477000000001
$ws.Range("M131").Value = -246697.5
$ws.Range("N131").Value = -16203.477

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3772.8333
$ws.Range("I80").Value = 4100.385
$ws.Range("J80").Value = 2921.2
$ws.Range("K80").Value = 4100.385
$ws.Range("L80").Value = 2921.2
$ws.Range("M80").Value = -3102.385
$ws.Range("N80").Value = -4917.2

$ws.Range("H83").Value = 3772.8333
$ws.Range("I83").Value = 4100.385
$ws.Range("J83").Value = 2921.2
$ws.Range("K83").Value = 20501.925
$ws.Range("L83").Value = 14606
$ws.Range("M83").Value = -15509.925
$ws.Range("N83").Value = -24590

$ws.Range("H132").Value = 2211.7104
$ws.Range("I132").Value = 1414.9584
$ws.Range("J132").Value = 3577.5715
$ws.Range("K132").Value = 4244.8752
$ws.Range("L132").Value = 10732.7145
$ws.Range("M132").Value = -1714.8752
$ws.Range("N132").Value = -15792.7145

$ws.Range("H138").Value = 22270.9
$ws.Range("J138").Value = 22270.9
$ws.Range("L138").Value = 22270.9
$ws.Range("N138").Value = -32550.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1425.4762
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 1528.1578
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 1528.1578
$ws.Range("M22").Value = -155
$ws.Range("N22").Value = -2118.1578

$ws.Range("H27").Value = 1425.4762
$ws.Range("I27").Value = 450
$ws.Range("J27").Value = 1528.1578
$ws.Range("K27").Value = 450
$ws.Range("L27").Value = 1528.1578
$ws.Range("M27").Value = -343
$ws.Range("N27").Value = -1742.1578

$ws.Range("H42").Value = 9000
$ws.Range("J42").Value = 9000
$ws.Range("L42").Value = 9000
$ws.Range("N42").Value = -10126

$ws.Range("H49").Value = 9000
$ws.Range("J49").Value = 9000
$ws.Range("L49").Value = 9000
$ws.Range("N49").Value = -9294

$ws.Range("H80").Value = 19875
$ws.Range("J80").Value = 19875
$ws.Range("L80").Value = 19875
$ws.Range("N80").Value = -22121

$ws.Range("H83").Value = 19875
$ws.Range("J83").Value = 19875
$ws.Range("L83").Value = 59625
$ws.Range("N83").Value = -70857

$ws.Range("H136").Value = 1943.6852
$ws.Range("I136").Value = 1222.2727
$ws.Range("J136").Value = 2439.6562
$ws.Range("K136").Value = 3666.8181
$ws.Range("L136").Value = 7318.9686
$ws.Range("M136").Value = -1116.8181
$ws.Range("N136").Value = -12418.9686

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 40027.25
$ws.Range("J138").Value = 40027.25
$ws.Range("L138").Value = 40027.25
$ws.Range("N138").Value = -50307.25
